# Live trading results update:
# Trade #54 opened (leadlag UP @ 21:32:05) and trades #24-#28 (leadlag,
# rows 20-24 on the "leadlag" sheet) closed via the 5-minute time exit.
# Propagate the newly-closed trades into "All Trades" and refresh the
# aggregate stats on "Summary" / "Comparison".

$wb  = $excel.ActiveWorkbook
$ll  = $wb.Worksheets.Item("leadlag")
$all = $wb.Worksheets.Item("All Trades")
$sum = $wb.Worksheets.Item("Summary")
$cmp = $wb.Worksheets.Item("Comparison")

# ---------------------------------------------------------------------
# 1) "leadlag" sheet: close out trades #24-#28 (rows 20-24)
# ---------------------------------------------------------------------

# Trade #24 (row 20)
$ll.Range("G20").Value = 67977.92436400001
$ll.Range("H20").Value = "CLOSED"
$ll.Range("I20").Value = 1.3869
$ll.Range("J20").Value = 13.87
$ll.Range("M20").Value = "time_exit_5min"
$ll.Range("N20").Value = 5

# Trade #25 (row 21)
$ll.Range("G21").Value = 68684.51604
$ll.Range("H21").Value = "CLOSED"
$ll.Range("I21").Value = 0.4877
$ll.Range("J21").Value = 4.88
$ll.Range("M21").Value = "time_exit_5min"
$ll.Range("N21").Value = 5

# Trade #26 (row 22)
$ll.Range("G22").Value = 68788.636335
$ll.Range("H22").Value = "CLOSED"
$ll.Range("I22").Value = -0.3256
$ll.Range("J22").Value = -3.26
$ll.Range("M22").Value = "time_exit_5min"
$ll.Range("N22").Value = 5

# Trade #27 (row 23)
$ll.Range("G23").Value = 69186.085081
$ll.Range("H23").Value = "CLOSED"
$ll.Range("I23").Value = -0.327
$ll.Range("J23").Value = -3.27
$ll.Range("M23").Value = "time_exit_5min"
$ll.Range("N23").Value = 5

# Trade #28 (row 24)
$ll.Range("G24").Value = 69071.66228800001
$ll.Range("H24").Value = "CLOSED"
$ll.Range("I24").Value = -0.238
$ll.Range("J24").Value = -2.38
$ll.Range("M24").Value = "time_exit_5min"
$ll.Range("N24").Value = 5

# ---------------------------------------------------------------------
# 2) "leadlag" sheet: append newly-opened trade #54 (row 44)
# ---------------------------------------------------------------------
$ll.Range("A44").Value = 54
$ll.Range("B44").Value = "'2026-02-16"
$ll.Range("C44").Value = "21:32:05"
$ll.Range("D44").Value = "leadlag"
$ll.Range("E44").Value = "UP"
$ll.Range("F44").Value = 68827.56
$ll.Range("H44").Value = "OPEN"
$ll.Range("I44").Value = 0
$ll.Range("J44").Value = 0
$ll.Range("K44").Value = 0.75
$ll.Range("L44").Value = "Binance leading with 0.124% move"
$ll.Range("N44").Value = 0

# ---------------------------------------------------------------------
# 3) "All Trades" sheet: mirror the 5 now-closed leadlag trades
#    (rows 20-24 on "leadlag") as new rows 25-29
# ---------------------------------------------------------------------
$ll.Range("A20:N20").Copy($all.Range("A25:N25"))
$ll.Range("A21:N21").Copy($all.Range("A26:N26"))
$ll.Range("A22:N22").Copy($all.Range("A27:N27"))
$ll.Range("A23:N23").Copy($all.Range("A28:N28"))
$ll.Range("A24:N24").Copy($all.Range("A29:N29"))

# ---------------------------------------------------------------------
# 4) "Summary" sheet: refresh OVERALL and leadlag aggregate rows
# ---------------------------------------------------------------------
$sum.Range("C2").Value = 28
$sum.Range("D2").Value = "'64.3%"
$sum.Range("E2").Value = "'+6.2783%"
$sum.Range("F2").Value = "'+0.2242%"

$sum.Range("C3").Value = 42
$sum.Range("D3").Value = "'33.3%"
$sum.Range("E3").Value = "'+4.5353%"
$sum.Range("F3").Value = "'+0.1080%"

# ---------------------------------------------------------------------
# 5) "Comparison" sheet: refresh leadlag row
# ---------------------------------------------------------------------
$cmp.Range("B2").Value = 42
$cmp.Range("C2").Value = "'33.3%"
$cmp.Range("D2").Value = "'2.42"
$cmp.Range("E2").Value = "'+0.5528%"
$cmp.Range("F2").Value = "'-0.3560%"
$cmp.Range("G2").Value = "'1.55"
